$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 values (B2:E2, G2) - F2 stays 0 (unchanged)
$ws.Range("B2").Value = 0.6545652718822623
$ws.Range("C2").Value = 1.626987699542094
$ws.Range("D2").Value = 0.7210945179870265
$ws.Range("E2").Value = 14773364.14517103
$ws.Range("G2").Value = 14773367.14781852

# Row 3 values (B3:E3, G3) - F3 stays 0 (unchanged)
$ws.Range("B3").Value = 0.003078177322033415
$ws.Range("C3").Value = 0.3048912486333797
$ws.Range("D3").Value = 3.223369029078222
$ws.Range("E3").Value = 2797.565817734744
$ws.Range("G3").Value = 2801.097156189777
